$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows to remove (1-indexed, matching the current sheet layout):
#   Row 5  -> 004751770 DILSON    60313.86
#   Row 7  -> 004308815 ZELI      22411.97
#   Row 9  -> 004479463 HENRIQUE   4411.39
#   Row 10 -> 004313254 GUSTAVO    4292
#   Row 11 -> 004332783 IRON       4028.3
#
# Delete bottom-to-top so earlier row numbers stay valid as we go.
$rowsToDelete = @(11, 10, 9, 7, 5)

foreach ($r in $rowsToDelete) {
    $ws.Rows.Item($r).Delete() | Out-Null
}
